$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = '29.910.45'
$ws.Cells.Item(2, 5).Value = '  +0.52%  '
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = '1.635.20'
$ws.Cells.Item(3, 5).Value = '  +1.01%  '
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = '1.00'
$ws.Cells.Item(4, 5).Value = '  +0.77%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '215.56'
$ws.Cells.Item(5, 5).Value = '  +1.34%  '
$ws.Cells.Item(6, 5).Value = '  +0.14%  '
$ws.Cells.Item(7, 5).Value = '  +0.83%  '
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '28.67'
$ws.Cells.Item(8, 5).Value = '  -1.75%  '
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.262'
$ws.Cells.Item(9, 5).Value = '  +1.22%  '
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '0.0611'
$ws.Cells.Item(10, 5).Value = '  +0.50%  '
$ws.Cells.Item(11, 5).Value = '  -1.04%  '
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '1.868.85'
$ws.Cells.Item(12, 5).Value = '  +0.86%  '
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '1.641.89'
$ws.Cells.Item(13, 5).Value = '  +1.10%  '
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '0.592'
$ws.Cells.Item(14, 5).Value = '  +4.61%  '
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '9.55'
$ws.Cells.Item(15, 5).Value = '  +6.89%  '
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '3.89'
$ws.Cells.Item(16, 5).Value = '  -0.25%  '
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '29.923.08'
$ws.Cells.Item(17, 5).Value = '  +0.52%  '
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '65.50'
$ws.Cells.Item(18, 5).Value = '  +1.81%  '
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '240.67'
$ws.Cells.Item(19, 5).Value = '  -0.36%  '
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '0.0₃0706'
$ws.Cells.Item(20, 5).Value = '  -0.38%  '
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '1.00'
$ws.Cells.Item(21, 5).Value = '  +0.69%  '
$ws.Cells.Item(22, 5).Value = '  +2.78%  '
$ws.Cells.Item(23, 5).Value = '  +1.45%  '
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '2.18'
$ws.Cells.Item(24, 5).Value = '  +3.27%  '
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '158.15'
$ws.Cells.Item(25, 5).Value = '  +1.95%  '
$ws.Cells.Item(26, 5).Value = '  -0.55%  '
$ws.Cells.Item(27, 5).Value = '  -1.07%  '
$ws.Cells.Item(28, 5).Value = '  +0.54%  '
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '1.00'
$ws.Cells.Item(29, 5).Value = '  +0.66%  '
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '0.0490'
$ws.Cells.Item(30, 5).Value = '  +0.43%  '
$ws.Cells.Item(31, 5).Value = '  +1.65%  '
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '3.39'
$ws.Cells.Item(32, 5).Value = '  +1.97%  '
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '3.20'
$ws.Cells.Item(33, 5).Value = '  -0.28%  '
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '1.424.89'
$ws.Cells.Item(34, 5).Value = '  +0.09%  '
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '1.70'
$ws.Cells.Item(35, 5).Value = '  +4.51%  '
$ws.Cells.Item(36, 5).Value = '  -0.93%  '
$ws.Cells.Item(37, 5).Value = '  -3.05%  '
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '0.0172'
$ws.Cells.Item(38, 5).Value = '  +1.40%  '
$ws.Cells.Item(39, 5).Value = '  +0.05%  '
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '76.27'
$ws.Cells.Item(40, 5).Value = '  +9.81%  '
$ws.Cells.Item(41, 5).Value = '  +0.15%  '
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '0.0503'
$ws.Cells.Item(42, 5).Value = '  -0.10%  '
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '0.834'
$ws.Cells.Item(43, 5).Value = '  +1.09%  '
$ws.Cells.Item(44, 5).Value = '  +1.02%  '
$ws.Cells.Item(45, 5).Value = '  +0.84%  '
$ws.Cells.Item(46, 5).Value = '  -0.68%  '
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '1.777.60'
$ws.Cells.Item(47, 5).Value = '  +0.86%  '
$ws.Cells.Item(48, 5).Value = '  -1.70%  '
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '48.74'
$ws.Cells.Item(49, 5).Value = '  -8.88%  '
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '92.77'
$ws.Cells.Item(50, 5).Value = '  +5.18%  '
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '0.0₆0111'
$ws.Cells.Item(51, 5).Value = '  +8.16%  '
